$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 (Euclid Managers): mark to pay and record this month's bill amount
$ws.Range("B5").Value = "x"
$ws.Range("D5").Value = 278.62

# Row 7 (Hernandez Lawn Service): mark to pay too (already flagged for mail/envelope)
$ws.Range("B7").Value = "x"

# Row 16 (UPS): this bill is settled - clear the pay mark, invoice number and amount
$ws.Range("B16").ClearContents()
$ws.Range("C16").ClearContents()
$ws.Range("D16").ClearContents()

# Leave the selection on the cell that was last edited
$ws.Range("B7").Select()
